$d = $word.ActiveDocument

# --- Edit 1: Replace the empty paragraph after "Insert a new instructor..." with 5 bold SQL paragraphs ---
$p1 = $d.Paragraphs(3)
$xml1 = '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">INSERT INTO instructor (INSTRUCTOR_ID, SALUTATION, FIRST_NAME, LAST_NAME, </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">           </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">STREET_ADDRESS, ZIP, CREATED_BY, CREATED_DATE, </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">          </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>MODIFIED_BY, MODIFIED_DATE)</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">VALUES (815, ''Mr.'', ''John'', ''Locke'', ''2342 Oceanic Way, Bayonne, NJ'', </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>07002, ''Ryan Buchanan'', SYSDATE, ''Ryan Buchanan'', SYSDATE);</w:t></w:r></w:p>'
$p1.Range.InsertXML($xml1)

# --- Edit 2: Add a left indent to the empty paragraph that follows "Create a new section ID..." ---
$p2 = $d.Paragraphs(10)
$p2.LeftIndent = 18

# --- Edit 3: styles.xml latent-style / DefaultParagraphFont touch-ups (done via the flat-OPC WordOpenXML roundtrip) ---
$xml = $d.WordOpenXML

$old1 = '<w:lsdException w:name="HTML Variable" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$new1 = '<w:lsdException w:name="HTML Variable" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Normal Table" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$xml = $xml.Replace($old1, $new1)

$old2 = '<w:lsdException w:name="Table Web 2" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$new2 = '<w:lsdException w:name="Table Web 2" w:semiHidden="1" w:unhideWhenUsed="1"/><w:lsdException w:name="Table Web 3" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$xml = $xml.Replace($old2, $new2)

$old3 = '<w:lsdException w:name="Table Grid" w:uiPriority="39"/>'
$new3 = '<w:lsdException w:name="Table Grid" w:uiPriority="39"/><w:lsdException w:name="Table Theme" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$xml = $xml.Replace($old3, $new3)

$old4 = '<w:style w:type="character" w:default="1" w:styleId="DefaultParagraphFont"><w:name w:val="Default Paragraph Font"/><w:uiPriority w:val="1"/><w:unhideWhenUsed/></w:style>'
$new4 = '<w:style w:type="character" w:default="1" w:styleId="DefaultParagraphFont"><w:name w:val="Default Paragraph Font"/><w:uiPriority w:val="1"/><w:semiHidden/><w:unhideWhenUsed/></w:style>'
$xml = $xml.Replace($old4, $new4)

$d.WordOpenXML = $xml

Write-Output "done"
